$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as plain text (e.g. "1.011", "28.006.85").
# Force text interpretation while assigning so Excel does not coerce these into numbers,
# then clear the temporary formatting so the cells keep their original (default) style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.006.85"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "1.884.48"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.78%  "
$ws.Range("D5").Value = "335.83"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("D7").Value = "0.4761"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("D8").Value = "0.3950"
$ws.Range("D9").Value = "46.94"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").Value = "0.08028"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "1.020"
$ws.Range("D12").Value = "21.93"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").Value = "1.886.32"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "6.065"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").Value = "7.207"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "1.014"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "88.49"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").Value = "0.06733"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "0.00001052"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "17.08"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "1.010"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "27.999.25"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("D23").Value = "5.510"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").Value = "11.00"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "2.350"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").Value = "2.122.18"
$ws.Range("E26").Value = "  +1.11%  "
$ws.Range("D27").Value = "159.15"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "19.93"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").Value = "2.110"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").Value = "5.519"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").Value = "121.67"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").Value = "0.9806"
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("D33").Value = "0.09575"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").Value = "3.640"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "5.348"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "1.364"
$ws.Range("E36").Value = "  -5.73%  "
$ws.Range("D37").Value = "0.06084"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "0.02254"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").Value = "1.210"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").Value = "8.220"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").Value = "0.5994"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "10.36"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").Value = "1.279"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").Value = "0.5674"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "12.26"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").Value = "3.353"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").Value = "0.06800"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").Value = "112.66"
$ws.Range("E51").Value = "  -1.55%  "

$ws.Range("D2:D51").ClearFormats()
